{"js": "// The data dictionary's \"Severity\" row description originally read:\n//   \"The severity of the data, ranging from 1 (least impact on traffic) to 4 (significant impact on traffic).\"\n// It should read \"...severity of the accident, ranging...\". The word\n// \"data\" is the only part that changes, so scope the search to the\n// specific table cell that holds the Severity description (the word\n// \"data\"/\"accident\" appears elsewhere in the document) and replace\n// just that one word, leaving the rest of the sentence untouched.\n\nconst table = context.document.body.tables.getFirst();\nconst cell = table.getCell(3, 3);\n\nconst results = cell.body.search(\"data\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target word to edit.\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\"accident\", \"Replace\");\nawait context.sync();\n", "ps1": "# The data dictionary's \"Severity\" row description originally read:\n#   \"The severity of the data, ranging from 1 (least impact on traffic) to 4 (significant impact on traffic).\"\n# It should read \"...severity of the accident, ranging...\". The word\n# \"data\" is the only part that changes, so scope the edit to the\n# specific table cell that holds the Severity description (the word\n# \"data\"/\"accident\" appears elsewhere in the document, e.g. \"Data Type\"\n# and \"Source of the raw accident data\") and replace just that one\n# word, leaving the rest of the sentence untouched.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$cell = $table.Cell(4, 4)\n$cellRange = $cell.Range\n\n$needle = \"data\"\n$cellText = $cellRange.Text\n$idx = $cellText.IndexOf($needle)\n\nif ($idx -lt 0) {\n    throw \"Could not find target word 'data' in the Severity description cell.\"\n}\n\n$wordStart = $cellRange.Start + $idx\n$wordEnd = $wordStart + $needle.Length\n$target = $d.Range($wordStart, $wordEnd)\n$target.Text = \"accident\"\n"}
